# Insert a new record row at row 427 (pushing the existing rows 427-512 down to 428-513)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 427:512 down by one row, same as Excel's Insert Row command.
$ws.Rows.Item(427).Insert()

# Fill in the new row 427 with the new record's data.
$ws.Range("A427").Value2 = 5
$ws.Range("B427").Value2 = "Macroferia Regional de Talca"
$ws.Range("C427").Value2 = "Maule"
$ws.Range("D427").Value2 = 45211
$ws.Range("E427").Value2 = 7
$ws.Range("F427").Value2 = 100112009
$ws.Range("G427").Value2 = "Acelga"
$ws.Range("H427").Value2 = "Sin especificar"
$ws.Range("I427").Value2 = "Primera"
$ws.Range("J427").Value2 = 400
$ws.Range("K427").Value2 = 1800
$ws.Range("L427").Value2 = 1800
$ws.Range("M427").Value2 = 1800
$ws.Range("N427").Value2 = "$/docena de atados (4 kilos)"
$ws.Range("O427").Value2 = "Región del Maule"
$ws.Range("P427").Value2 = 450
$ws.Range("Q427").Value2 = 4
$ws.Range("R427").Value2 = "Hortaliza"
